# Form the consolidated report: fill in the "Absent" (column H) values
# based on the "Real" (column E) attendance flag for each data row.
# Absent = 1 when Real == 0, otherwise Absent = 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 3; $r -le $lastRow; $r++) {
    $real = $ws.Cells.Item($r, 5).Value2
    if ($real -eq 0) {
        $ws.Cells.Item($r, 8).Value2 = 1
    } else {
        $ws.Cells.Item($r, 8).Value2 = 0
    }
}
